$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a value as text, preserving the cell's original (unstyled)
# appearance. Excel's COM layer auto-detects numeric-looking strings (like
# "1.00" or "0.999") and would otherwise coerce them into numbers, losing
# the exact textual representation used in the source data. Prefixing with
# a quote forces text entry; the cell Style is saved/restored so no
# quote-prefix formatting flag leaks into the saved style table.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "52.081.23"
$ws.Range("E2").Value = "  +0.73%  "
Set-TextValue $ws.Range("D3") "2.872.05"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "350.22"
$ws.Range("E5").Value = "  -0.84%  "
Set-TextValue $ws.Range("D6") "112.45"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("E7").Value = "  +1.09%  "
Set-TextValue $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  -0.01%  "
Set-TextValue $ws.Range("D9") "0.619"
$ws.Range("E9").Value = "  +1.87%  "
Set-TextValue $ws.Range("D10") "40.22"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").Value = "  -0.63%  "
Set-TextValue $ws.Range("D12") "0.0850"
$ws.Range("E12").Value = "  +1.52%  "
Set-TextValue $ws.Range("D13") "20.05"
$ws.Range("E13").Value = "  +0.14%  "
Set-TextValue $ws.Range("D14") "7.85"
$ws.Range("E14").Value = "  +2.02%  "
Set-TextValue $ws.Range("D15") "3.325.31"
$ws.Range("E15").Value = "  +3.14%  "
Set-TextValue $ws.Range("D16") "0.994"
$ws.Range("E16").Value = "  +6.99%  "
Set-TextValue $ws.Range("D17") "2.879.28"
$ws.Range("E17").Value = "  +2.34%  "
Set-TextValue $ws.Range("D18") "52.070.73"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("E20").Value = "  -1.51%  "
Set-TextValue $ws.Range("D21") "13.59"
$ws.Range("E22").Value = "  +0.99%  "
Set-TextValue $ws.Range("D23") "70.90"
$ws.Range("E23").Value = "  +1.37%  "
Set-TextValue $ws.Range("D24") "270.18"
$ws.Range("E25").Value = "  +1.59%  "
Set-TextValue $ws.Range("D26") "26.53"
$ws.Range("E26").Value = "  +1.73%  "
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("E28").Value = "  -1.35%  "
Set-TextValue $ws.Range("D29") "10.57"
$ws.Range("E29").Value = "  +2.89%  "
Set-TextValue $ws.Range("D30") "38.83"
$ws.Range("E30").Value = "  +3.98%  "
Set-TextValue $ws.Range("D31") "6.26"
$ws.Range("E31").Value = "  +0.79%  "
Set-TextValue $ws.Range("D32") "52.34"
$ws.Range("E32").Value = "  +1.13%  "
Set-TextValue $ws.Range("D33") "5.84"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  -0.41%  "
Set-TextValue $ws.Range("D35") "0.0897"
$ws.Range("E35").Value = "  +7.46%  "
Set-TextValue $ws.Range("D36") "0.999"
$ws.Range("E36").Value = "  -0.07%  "
Set-TextValue $ws.Range("D37") "1.89"
$ws.Range("E37").Value = "  -15.41%  "
$ws.Range("E38").Value = "  +5.96%  "
Set-TextValue $ws.Range("D39") "18.70"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("E40").Value = "  +3.25%  "
Set-TextValue $ws.Range("D41") "2.64"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("E42").Value = "  +1.39%  "
Set-TextValue $ws.Range("D43") "121.92"
$ws.Range("E43").Value = "  +1.19%  "
Set-TextValue $ws.Range("D44") "22.43"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("E46").Value = "  +5.04%  "
Set-TextValue $ws.Range("D47") "2.171.02"
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("E48").Value = "  +5.82%  "
$ws.Range("E49").Value = "  +11.33%  "
Set-TextValue $ws.Range("D50") "0.961"
$ws.Range("E50").Value = "  +5.77%  "
Set-TextValue $ws.Range("D51") "0.0322"
$ws.Range("E51").Value = "  +12.90%  "
